$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected Price (D) and Volume (E) columns keep their original
# text formatting (e.g. "1.00", "0.999", leading/trailing spaces) instead of
# being auto-converted to numbers by Excel when the new value looks numeric.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.419.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.071.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.070.81'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.439'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.601.89'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.29'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.45%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.485.53'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.067.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.60'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '338.32'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.80'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.169'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.32'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.79'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.33'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.03%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.11'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.76'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.85%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.78%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.680'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.111.14'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.71'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.267.31'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.984'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.04%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.42'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.48%  '
